# Source titles are added (rows removed: Creating a Dashboard, Formatting data,
# Uploading data, Creating a story) -- delete rows 27-30 on the topic_ids sheet,
# shifting the remaining rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("topic_ids")

$ws.Rows("27:30").Delete()
